$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.308.62'
$ws.Range('E2').Value = '  -2.93%  '
$ws.Range('D3').Value = '3.763.14'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.79'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.66'
$ws.Range('E6').Value = '  -3.66%  '
$ws.Range('D7').Value = '3.761.70'
$ws.Range('E7').Value = '  -1.41%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('E10').Value = '  -3.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.37'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -5.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.82'
$ws.Range('E14').Value = '  -3.24%  '
$ws.Range('D15').Value = '4.398.59'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').Value = '3.767.74'
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').Value = '67.305.61'
$ws.Range('E17').Value = '  -2.84%  '
$ws.Range('E18').Value = '  -3.25%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.22'
$ws.Range('E21').Value = '  -8.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '455.82'
$ws.Range('E22').Value = '  -3.76%  '
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.03'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('E27').Value = '  -6.62%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.98'
$ws.Range('E29').Value = '  -3.21%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.65'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.18'
$ws.Range('E32').Value = '  -4.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  -4.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.14'
$ws.Range('E34').Value = '  -3.08%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '3.719.46'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0995'
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('E38').Value = '  -8.59%  '
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.991'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  -3.31%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.67'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('E45').Value = '  -4.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.77'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.33'
$ws.Range('E47').Value = '  -3.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '147.14'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('E49').Value = '  -8.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '388.59'
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('D51').Value = '2.742.55'
$ws.Range('E51').Value = '  +1.41%  '
